# Replace the old "Mask" icon row on the "Main Icons" sheet with a new
# "Eye" icon row, and add a brand-new "Hidden Eye" icon row right after it.
#
# Cell-write order matters here: it determines the order new entries land
# in the shared-strings table, so we write B33 (the URL) before A33 (the
# label) for row 33 -- mirroring how the row above it (URL before label)
# was originally authored -- and then A34 before B34 (label before URL,
# the "normal" order used everywhere else in the sheet) for the new row.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Main Icons")

$ws.Range("B33").Value = "http://www.flaticon.com/free-icon/visible-opened-eye-interface-option_58976"
$ws.Range("A33").Value = "Eye"

$ws.Range("A34").Value = "Hidden Eye"
$ws.Range("B34").Value = "http://www.flaticon.com/free-icon/invisible_59394"

# Make "Main Icons" the active sheet/tab again (it was "Videos" before),
# with B34 -- the last cell we touched -- selected.
$ws.Activate()
$ws.Range("B34").Select()
